$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed values (B,C,D,E,G) for rows 2-8. Column F is unchanged.
$data = @{
    2 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 14.05633640148523)
    3 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    4 = @(0.6606524410359556, 250555.8564151394, 0.7527432677738641, 10.19245300693656, 250567.4622638551)
    5 = @(0.0006408296065709695, 0.002571899574220771, 3.537761648806719, 1133.036916526867, 1136.577890904855)
    6 = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086)
    7 = @(0.1190320826869504, 10.34677158129881, 0.1494219747398047, 1133.036916526867, 1143.652142165593)
    8 = @(0.04271373187048222, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.596514762964814)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}

$wb.Save()
